$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 11 - "listado de articulos habilitados" (GET articulosHabilitados)
# ---------------------------------------------------------------------------
# Seed formatting by copying the most similar existing row (row 10 matches the
# target style pattern for columns A-F exactly), then overwrite contents.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("B10").Copy($ws.Range("B11"))
$ws.Range("C10").Copy($ws.Range("C11"))
$ws.Range("D10").Copy($ws.Range("D11"))
$ws.Range("E10").Copy($ws.Range("E11"))
$ws.Range("F10").Copy($ws.Range("F11"))

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "listado de articulos habilitados"
$ws.Range("C11").Value = "http://localhost:8089/b-salesforce/rest/articulosHabilitados"
# D11 keeps the inherited "GET" text (shared with D10)
$f11 = @"
{
    "success": true,
    "message": "Consulta Exitosa.",
    "result": [
        {
            "id": 1,
            "descripcionArticulo": "ANILLO CUADRADO",
            "codigoArticulo": "JOAN0.0.1",
            "metodoCosto": "95.0",
            "precio": 100,
            "precioCosto": null,
            "upc": null,
            "nivelReorden": 0,
            "cantidadReorden": 3,
            "nSerie": 0,
            "fotografia": null,
            "fechaDesde": 1410386400000,
            "fechaHasta": null,
            "usuarioAct": "JSON"
        },
        {
            "id": 2,
            "descripcionArticulo": "ANILLO RECTANGULAR",
            "codigoArticulo": "JOAN0.0.2",
            "metodoCosto": "250.0",
            "precio": 15000,
            "precioCosto": null,
            "upc": null,
            "nivelReorden": 0,
            "cantidadReorden": 3,
            "nSerie": 0,
            "fotografia": null,
            "fechaDesde": 1410386400000,
            "fechaHasta": null,
            "usuarioAct": "JSON"
        },
        {
            "id": 3,
            "descripcionArticulo": "ANILLO ZEBRA",
            "codigoArticulo": "JOAN0.0.3",
            "metodoCosto": "500.0",
            "precio": 2500,
            "precioCosto": null,
            "upc": null,
            "nivelReorden": 0,
            "cantidadReorden": 3,
            "nSerie": 0,
            "fotografia": null,
            "fechaDesde": 1410386400000,
            "fechaHasta": null,
            "usuarioAct": "JSON"
        }
    ]
}
"@
$ws.Range("F11").Value = $f11
$ws.Rows.Item(11).RowHeight = 409.5

# ---------------------------------------------------------------------------
# New row 12 - "Listado de proveedores habilitados" (GET proveedorMovimiento)
# ---------------------------------------------------------------------------
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("B8").Copy($ws.Range("B12"))
$ws.Range("C8").Copy($ws.Range("C12"))
$ws.Range("D8").Copy($ws.Range("D12"))
$ws.Range("E8").Copy($ws.Range("E12"))
$ws.Range("F8").Copy($ws.Range("F12"))

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Listado de proveedores habilitados"
$ws.Range("C12").Value = "http://localhost:8089/b-salesforce/rest/proveedorMovimiento"
# D12 keeps the inherited "GET" text (shared with D8)
$f12 = @"
{
    "success": true,
    "message": "Consulta exitosa.",
    "result": [
        {
            "id": 1,
            "nombre": "COMPANEX SRL",
            "numeroDocumento": "4851800011"
        },
        {
            "id": 2,
            "nombre": "YANAPAX SRL",
            "numeroDocumento": "48484850001"
        },
        {
            "id": 3,
            "nombre": "infoglobal",
            "numeroDocumento": "78787878778"
        },
        {
            "id": 4,
            "nombre": "ENATEX",
            "numeroDocumento": "4851800"
        }
    ]
}
"@
$ws.Range("F12").Value = $f12
$ws.Rows.Item(12).RowHeight = 390

# Wire up the hyperlink for C12, then restore the original hyperlink-cell
# formatting (Hyperlinks.Add stamps its own style on the target cell).
$ws.Hyperlinks.Add($ws.Range("C12"), "http://localhost:8089/b-salesforce/rest/proveedorMovimiento") | Out-Null
$ws.Range("C8").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Viewport: scroll down and move the selection onto the newly added data.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("E11").Select()
